$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (team specific matrix VMI_A)
# Row 2
$ws.Range("B2").Value = 0.2210526315789474
$ws.Range("C2").Value = 0.5228070175438596
$ws.Range("J2").Value = 0.007017543859649123
$ws.Range("P2").Value = 0.1508771929824561
$ws.Range("S2").Value = 0.09824561403508772

# Row 3
$ws.Range("B3").Value = 0.006756756756756757
$ws.Range("C3").Value = 0.01351351351351351
$ws.Range("J3").Value = 0.03378378378378379
$ws.Range("P3").Value = 0.7702702702702703
$ws.Range("S3").Value = 0.1756756756756757

# Row 4
$ws.Range("J4").Value = 0.08571428571428572
$ws.Range("P4").Value = 0.6285714285714286
$ws.Range("S4").Value = 0.2857142857142857

# Row 5
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5

# Row 6
$ws.Range("B6").Value = 0.03389830508474576
$ws.Range("D6").Value = 0.01271186440677966
$ws.Range("F6").Value = 0.05508474576271186
$ws.Range("J6").Value = 0.2754237288135593
$ws.Range("O6").Value = 0.0211864406779661
$ws.Range("Q6").Value = 0.1483050847457627
$ws.Range("R6").Value = 0.05508474576271186
$ws.Range("S6").Value = 0.3983050847457627

# Row 7
$ws.Range("B7").Value = 0.08247422680412371
$ws.Range("D7").Value = 0.02577319587628866
$ws.Range("E7").Value = 0.005154639175257732
$ws.Range("F7").Value = 0.04123711340206185
$ws.Range("J7").Value = 0.1185567010309278
$ws.Range("O7").Value = 0.01030927835051546
$ws.Range("Q7").Value = 0.1597938144329897
$ws.Range("R7").Value = 0.04123711340206185
$ws.Range("S7").Value = 0.5154639175257731

# Row 8
$ws.Range("B8").Value = 0.06349206349206349
$ws.Range("D8").Value = 0.007054673721340388
$ws.Range("F8").Value = 0.04761904761904762
$ws.Range("J8").Value = 0.1005291005291005
$ws.Range("O8").Value = 0.02292768959435626
$ws.Range("Q8").Value = 0.1798941798941799
$ws.Range("R8").Value = 0.1111111111111111
$ws.Range("S8").Value = 0.4673721340388007

# Row 9
$ws.Range("B9").Value = 0.08118081180811808
$ws.Range("D9").Value = 0.01107011070110701
$ws.Range("E9").Value = 0.007380073800738007
$ws.Range("F9").Value = 0.06273062730627306
$ws.Range("J9").Value = 0.07749077490774908
$ws.Range("O9").Value = 0.01107011070110701
$ws.Range("Q9").Value = 0.1881918819188192
$ws.Range("R9").Value = 0.07749077490774908
$ws.Range("S9").Value = 0.4833948339483395

# Row 10
$ws.Range("B10").Value = 0.09900990099009901
$ws.Range("D10").Value = 0.01768033946251768
$ws.Range("E10").Value = 0.0007072135785007072
$ws.Range("F10").Value = 0.07001414427157002
$ws.Range("J10").Value = 0.1004243281471004
$ws.Range("O10").Value = 0.02121640735502122
$ws.Range("Q10").Value = 0.2072135785007072
$ws.Range("R10").Value = 0.06577086280056577
$ws.Range("S10").Value = 0.417963224893918

# Row 11
$ws.Range("G11").Value = 0.1597633136094675
$ws.Range("J11").Value = 0.121301775147929
$ws.Range("K11").Value = 0.2307692307692308
$ws.Range("L11").Value = 0.4585798816568047
$ws.Range("S11").Value = 0.02958579881656805

# Row 12
$ws.Range("G12").Value = 0.6890243902439024
$ws.Range("J12").Value = 0.2073170731707317
$ws.Range("K12").Value = 0.006097560975609756
$ws.Range("L12").Value = 0.03658536585365853
$ws.Range("S12").Value = 0.06097560975609756

# Row 13
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.08333333333333333

# Row 15
$ws.Range("F15").Value = 0.0228310502283105
$ws.Range("H15").Value = 0.1643835616438356
$ws.Range("I15").Value = 0.0730593607305936
$ws.Range("J15").Value = 0.3424657534246575
$ws.Range("K15").Value = 0.0365296803652968
$ws.Range("M15").Value = 0.0136986301369863
$ws.Range("O15").Value = 0.0867579908675799
$ws.Range("S15").Value = 0.2602739726027397

# Row 16
$ws.Range("F16").Value = 0.03428571428571429
$ws.Range("H16").Value = 0.2114285714285714
$ws.Range("I16").Value = 0.1085714285714286
$ws.Range("J16").Value = 0.3885714285714286
$ws.Range("K16").Value = 0.1085714285714286
$ws.Range("M16").Value = 0.01142857142857143
$ws.Range("O16").Value = 0.01142857142857143
$ws.Range("S16").Value = 0.1257142857142857

# Row 17
$ws.Range("F17").Value = 0.01377952755905512
$ws.Range("H17").Value = 0.1909448818897638
$ws.Range("I17").Value = 0.1023622047244094
$ws.Range("J17").Value = 0.4074803149606299
$ws.Range("K17").Value = 0.08858267716535433
$ws.Range("M17").Value = 0.00984251968503937
$ws.Range("N17").Value = 0.001968503937007874
$ws.Range("O17").Value = 0.05511811023622047
$ws.Range("S17").Value = 0.1299212598425197

# Row 18
$ws.Range("F18").Value = 0.02512562814070352
$ws.Range("H18").Value = 0.2060301507537688
$ws.Range("I18").Value = 0.1055276381909548
$ws.Range("J18").Value = 0.3668341708542713
$ws.Range("K18").Value = 0.09547738693467336
$ws.Range("M18").Value = 0.01507537688442211
$ws.Range("O18").Value = 0.06532663316582915
$ws.Range("S18").Value = 0.1206030150753769

# Row 19
$ws.Range("F19").Value = 0.01292307692307692
$ws.Range("H19").Value = 0.2258461538461538
$ws.Range("I19").Value = 0.1009230769230769
$ws.Range("J19").Value = 0.3735384615384615
$ws.Range("K19").Value = 0.09661538461538462
$ws.Range("M19").Value = 0.02215384615384615
$ws.Range("N19").Value = 0.0006153846153846154
$ws.Range("O19").Value = 0.04615384615384616
$ws.Range("S19").Value = 0.1212307692307692

